# Update Brake Service (FSG Rules), Update Pins
#
# Renames pin-function labels and swaps the old Arduino-style pin names
# (A0-A5, D3/D5/D6/D7/D9/D10) with the new STM32-style pin names
# (PA_x / PB_x / PC_x / ... "(links)" / "(rechts)") on the "Pinbelegung"
# worksheet (Tabelle1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- RPM sensor rows (Pedalerie block) - pins were previously empty ---
$ws.Range("D22").Value = "PB_2 (rechts)"
$ws.Range("D23").Value = "PB_4 (rechts)"

# --- Pedalerie block (rows 16-21) ---
$ws.Range("D17").Value = "PA_0 (links)"
$ws.Range("D18").Value = "PB_0 (links)"
$ws.Range("D19").Value = "PC_1 (links)"
$ws.Range("D16").Value = "PB_1 (rechts)"
$ws.Range("D20").Value = "PA_3 (rechts)"
$ws.Range("D21").Value = "PC_4 (rechts)"

# --- Taster / LED block (rows 6-10) ---
$ws.Range("C6").Value = "Taster Start"
$ws.Range("C7").Value = "Taster Reset"
$ws.Range("C8").Value = "LED Grün"
$ws.Range("C9").Value = "LED Gelb"
$ws.Range("C10").Value = "LED Rot"
$ws.Range("D10").Value = "PB_4 (rechts)"
$ws.Range("D9").Value = "PB_5 (rechts)"
$ws.Range("D8").Value = "PA_10 (rechts)"
$ws.Range("D7").Value = "PB_10 (rechts)"
$ws.Range("D6").Value = "PC_7 (rechts)"

# --- LV-Box block (rows 28-33) ---
$ws.Range("D28").Value = "PG_4 (rechts)"
$ws.Range("D29").Value = "PE_13 (rechts)"
$ws.Range("D30").Value = "PD_10 (rechts)"
$ws.Range("D31").Value = "PE_10 (rechts)"
$ws.Range("D32").Value = "PF_4 (rechts)"
$ws.Range("D33").Value = "PB_1 (rechts)"

# --- Update the view/selection state to match the author's last cursor position ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D34").Select() | Out-Null
